$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Poste" (C), "Departement" (D) text values, and "Salaire De Base" (J) numeric values ---
# Row => Poste, Departement, Salaire
$data = @{
    2  = @("pas cadre", "Mecanicien", 25)
    3  = @("pas cadre", "agent de recrutement", 27)
    4  = @("cadre", "Logisticien", 71)
    5  = @(" cadre", "Agent d'audit", 55)
    6  = @("cadre", "Assistant de direction", 45)
    7  = @("pas cadre", "Developpeur", 72)
    8  = @("pas cadre", "Mecanicien", 20)
    9  = @("pas cadre", "Mecanicien", 17)
    10 = @("pas cadre", "Agent", 63)
    11 = @("pas cadre", "sécretaire", 18)
    12 = @("pas cadre", "Mecanicien", 14)
    13 = @("pas cadre", "Community manager", 13)
    14 = @("pas cadre", "Mecanicien", 25)
    15 = @("pas cadre", "Mecanicien", 17)
    16 = @("pas cadre", "Mecanicien", 30)
    17 = @("pas cadre", "Manuttentionaire", 37)
    18 = @("pas cadre", "Agent de recouvrement", 16)
    19 = @("cadre", "Informaticien", 77)
    20 = @("pas cadre", "agent d entretien", 34)
    21 = @("cadre", "gestionnaire RH", 80)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("D$row").Value = $vals[1]
    $ws.Range("J$row").Value = $vals[2]
}

# --- Adjust column E (Ville) width to 8 characters (stored width) ---
# Excel's ColumnWidth (chars) maps to the stored <col width> via +5/6,
# so subtract that padding to land exactly on a stored width of 8.
$ws.Columns.Item(5).ColumnWidth = 8 - 5/6

# --- Update sheet view: scroll the window so column B is the leftmost visible
#     column, and move the selection to J3 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("J3").Select()
